$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2906.4167
$ws.Range("J17").Value = 1999.8334
$ws.Range("L17").Value = 5999.5002
$ws.Range("N17").Value = -6335.5002

$ws.Range("H33").Value = 337.1
$ws.Range("I33").Value = 103.5
$ws.Range("K33").Value = 103.5
$ws.Range("M33").Value = 125.5

$ws.Range("H116").Value = 5860.304
$ws.Range("I116").Value = 5451.5
$ws.Range("K116").Value = 5451.5
$ws.Range("M116").Value = -2009.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1400
$ws.Range("I61").Value = 1400
$ws.Range("K61").Value = 1400
$ws.Range("M61").Value = -1188

$ws.Range("H97").Value = 665
$ws.Range("I97").Value = 665
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 665
$ws.Range("L97").ClearContents()
$ws.Range("N97").Value = 0
$ws.Range("M97").Value = -169

$ws.Range("H102").Value = 923.44446
$ws.Range("J102").Value = 1352.7273
$ws.Range("L102").Value = 1352.7273
$ws.Range("N102").Value = -4596.7273

$ws.Range("H132").Value = 4903
$ws.Range("I132").Value = 1806
$ws.Range("K132").Value = 5418
$ws.Range("M132").Value = -2888

$ws.Range("H136").Value = 1400
$ws.Range("I136").Value = 1400
$ws.Range("K136").Value = 4200
$ws.Range("M136").Value = -1650

$ws.Range("H140").Value = 107397.8
$ws.Range("J140").Value = 107397.8
$ws.Range("L140").Value = 107397.8
$ws.Range("N140").Value = -117757.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3248.0588
$ws.Range("I134").Value = 2944.1428
$ws.Range("K134").Value = 8832.428400000001
$ws.Range("M134").Value = -6297.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 33799
$ws.Range("J50").Value = 33799
$ws.Range("L50").Value = 33799
$ws.Range("N50").Value = -35049

$ws.Range("H51").Value = 46666.332
$ws.Range("J51").Value = 46666.332
$ws.Range("L51").Value = 46666.332
$ws.Range("N51").Value = -48138.332

$ws.Range("H58").Value = 3718.5881
$ws.Range("I58").Value = 2801
$ws.Range("K58").Value = 2801
$ws.Range("M58").Value = -2598

$ws.Range("H60").Value = 26142.572
$ws.Range("J60").Value = 49999
$ws.Range("L60").Value = 49999
$ws.Range("N60").Value = -51021

$ws.Range("H61").Value = 46666.332
$ws.Range("J61").Value = 46666.332
$ws.Range("L61").Value = 46666.332
$ws.Range("N61").Value = -47362.332

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0

$ws.Range("H105").Value = 2768.1428
$ws.Range("I105").Value = 836.7143
$ws.Range("J105").Value = 4699.5713
$ws.Range("K105").Value = 836.7143
$ws.Range("L105").Value = 4699.5713
$ws.Range("M105").Value = 910.2857
$ws.Range("N105").Value = -8193.5713

$ws.Range("H132").Value = 4640.7856
$ws.Range("I132").Value = 2902.625
$ws.Range("K132").Value = 8707.875
$ws.Range("M132").Value = -6177.875

$ws.Range("H134").Value = 3541.4783
$ws.Range("I134").Value = 3049.5625
$ws.Range("K134").Value = 9148.6875
$ws.Range("M134").Value = -6613.6875

$ws.Range("H136").Value = 3718.5881
$ws.Range("I136").Value = 2801
$ws.Range("K136").Value = 8403
$ws.Range("M136").Value = -5853

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2699.8333
$ws.Range("I62").Value = 1399.5
$ws.Range("J62").Value = 3350
$ws.Range("K62").Value = 4198.5
$ws.Range("L62").Value = 10050
$ws.Range("M62").Value = -3512.5
$ws.Range("N62").Value = -11422

$ws.Range("H65").Value = 2699.8333
$ws.Range("I65").Value = 1399.5
$ws.Range("J65").Value = 3350
$ws.Range("K65").Value = 12595.5
$ws.Range("L65").Value = 30150
$ws.Range("M65").Value = -9163.5
$ws.Range("N65").Value = -37014

$ws.Range("H113").Value = 1160.7059
$ws.Range("J113").Value = 990.3333
$ws.Range("L113").Value = 2970.9999
$ws.Range("N113").Value = -7310.9999

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = 0

$ws.Range("H137").Value = 9606.6
$ws.Range("I137").Value = 9266.666999999999
$ws.Range("K137").Value = 27800.001
$ws.Range("M137").Value = -22700.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 919.25
$ws.Range("I2").Value = 228.5
$ws.Range("K2").Value = 228.5
$ws.Range("M2").Value = -115.5

$ws.Range("H35").Value = 5333333.5
$ws.Range("J35").Value = 5000000
$ws.Range("L35").Value = 5000000
$ws.Range("N35").Value = -5000596

$ws.Range("H132").Value = 4557.6
$ws.Range("I132").Value = 2950
$ws.Range("K132").Value = 8850
$ws.Range("M132").Value = -6320

$ws.Range("H136").Value = 23297.938
$ws.Range("J136").Value = 23297.938
$ws.Range("L136").Value = 69893.814
$ws.Range("N136").Value = -74993.814

$ws.Range("H141").Value = 61716.332
$ws.Range("J141").Value = 61716.332
$ws.Range("L141").Value = 61716.332
$ws.Range("N141").Value = -72076.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1758
$ws.Range("I30").Value = 1016
$ws.Range("K30").Value = 1016
$ws.Range("M30").Value = -908

$ws.Range("H46").Value = 3642
$ws.Range("I46").Value = 2623.625
$ws.Range("K46").Value = 2623.625
$ws.Range("M46").Value = -2435.625

$ws.Range("H104").Value = 30684.75
$ws.Range("J104").Value = 30684.75
$ws.Range("L104").Value = 30684.75
$ws.Range("N104").Value = -37672.75

$ws.Range("H122").Value = 5043.4443
$ws.Range("I122").Value = 3397
$ws.Range("J122").Value = 5866.6665
$ws.Range("K122").Value = 10191
$ws.Range("L122").Value = 17599.9995
$ws.Range("M122").Value = -7741
$ws.Range("N122").Value = -22499.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 29989.5
$ws.Range("I3").Value = 29989
$ws.Range("J3").Value = 29990
$ws.Range("K3").Value = 29989
$ws.Range("L3").Value = 29990
$ws.Range("M3").Value = -29875
$ws.Range("N3").Value = -30218

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0

$ws.Range("H81").Value = 1399.8
$ws.Range("I81").Value = 1399.8
$ws.Range("K81").Value = 2799.6
$ws.Range("M81").Value = -1738.6

$ws.Range("H84").Value = 1399.8
$ws.Range("I84").Value = 1399.8
$ws.Range("K84").Value = 13998
$ws.Range("M84").Value = -8694

$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -43494

$ws.Range("H122").Value = 705.86664
$ws.Range("I122").Value = 633.25
$ws.Range("K122").Value = 1899.75
$ws.Range("M122").Value = 550.25

$ws.Range("H126").Value = 2026.0834
$ws.Range("I126").Value = 814.1111
$ws.Range("K126").Value = 2442.3333
$ws.Range("M126").Value = 27.66670000000022

$ws.Range("H132").Value = 1598.5454
$ws.Range("I132").Value = 1598.5454
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4795.6362
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -2265.6362

$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360
